$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(14, 'KCH', 'Bakalářská práce', 'Bakalářská práce'),
    @(15, 'KFY', 'Bakalářská práce', 'Bakalářská práce'),
    @(16, 'KMA', 'Bakalářská práce', 'Bakalářská práce'),
    @(17, 'KFY', 'Bakalářská práce I', 'Bakalářská práce I'),
    @(18, 'CNB', 'Bakalářská práce I', 'Bakalářská práce I'),
    @(19, 'KFY', 'Bakalářská práce II', 'Bakalářská práce II'),
    @(20, 'CNB', 'Bakalářská práce II', 'Bakalářská práce II'),
    @(28, 'KI', 'Diplomová práce', 'Diplomová práce'),
    @(29, 'KMA', 'Diplomová práce', 'Diplomová práce'),
    @(41, 'KFY', 'Bakalářská práce z fyziky', 'Bakalářská práce z fyziky'),
    @(42, 'KCH', 'Bakalářská práce z chemie', 'Bakalářská práce z chemie'),
    @(48, 'KFY', 'PVK - Automatizované měřící systémy', 'Automatizované měřící systémy'),
    @(49, 'KCH', 'Biochemická cvičení', 'Biochemická cvičení'),
    @(50, 'KCH', 'Diplomová práce z chemie', 'Diplomová práce z chemie'),
    @(51, 'KFY', 'Diplomová práce fyziky', 'Diplomová práce'),
    @(79, 'KFY', 'Bakalářská práce z fyziky', 'Bakalářská práce z fyziky'),
    @(80, 'KBI', 'Bakalářská práce z biologie', 'Bakalářská práce'),
    @(84, 'KPRF', 'Anglický jazyk pro doktorandy', 'Anglický jazyk pro doktorandy'),
    @(85, 'KFY', 'Anglický jazyk pro doktorandy', 'Anglický jazyk pro doktorandy'),
    @(86, 'KBI', 'SZZ - biologie s didaktikou', 'Biologie s didaktikou'),
    @(87, 'KGEO', 'SZZ - Geografie s didaktikou', 'Geografie s didaktikou pro ZŠ'),
    @(88, 'KMA', 'SZZ - matematika s didaktikou', 'Matematika s didaktikou'),
    @(90, 'KMA', 'SZZ - matematika s didaktikou', 'Matematika s didaktikou'),
    @(91, 'KGEO', 'SZZ - Geografie s didaktikou', 'Geografie s didaktikou pro SŠ'),
    @(94, 'KFY', 'Aplikovaná fyzika', 'Aplikovaná fyzika'),
    @(95, 'KGEO', 'Aplikovaná geografie - Cestovní ruch', 'Aplikovaná geografie - Cestovní ruch'),
    @(96, 'KFY', 'Nanotechnologie', 'Nanotechnologie'),
    @(97, 'KI', 'Informační technologie', 'Informační technologie'),
    @(98, 'KGEO', 'Aplikovaná geografie - Krajina a GIS', 'Aplikovaná geografie - Krajina a GIS'),
    @(100, 'KGEO', 'SZZ - Evropská integrace a EU', 'Evropská integrace a Evropská unie'),
    @(101, 'KFY', 'Počítačové modelování', 'Počítačové modelování'),
    @(104, 'KBI', 'Biologie s didaktikou pro SŠ', 'Biologie s didaktikou pro SŠ'),
    @(106, 'KGEO', 'Geografie krajiny a GIS', 'Geografie krajiny a GIS'),
    @(110, 'KFY', 'SZZ - numerická matematika', 'Numerická matematika'),
    @(111, 'KGEO', 'SZZ - Reg. geografie a reg. rozvoj Česka', 'Regionální geografie a regionální rozvoj Česka'),
    @(113, 'KGEO', 'Geografie pro vzdělávání', 'Geografie pro vzdělávání'),
    @(114, 'KCH', 'Chemie pro vzdělávání', 'Chemie pro vzdělávání'),
    @(115, 'KMA', 'Matematika pro vzdělávání', 'Matematika pro vzdělávání'),
    @(116, 'KBI', 'Biologie pro vzdělávání', 'Biologie pro vzdělávání'),
    @(117, 'KFY', 'Fyzika pro vzdělávání', 'Fyzika pro vzdělávání'),
    @(118, 'KCH', 'Analytická chemie', 'Analytická chemie'),
    @(119, 'KBI', 'SZZ - Biologie', 'Biologie'),
    @(120, 'KFY', 'Fyzika', 'Fyzika'),
    @(121, 'KFY', 'Fyzika', 'Fyzika'),
    @(122, 'KGEO', 'SZZ- Regionální geografie Evropy a světa', 'Regionální geografie Evropy a světa'),
    @(123, 'KCH', 'Toxikologie', 'Toxikologie'),
    @(124, 'KMA', 'Matematika s didaktikou pro ZŠ', 'Matematika s didaktikou pro ZŠ'),
    @(127, 'KFY', 'Fyzika', 'Fyzika'),
    @(128, 'KBI', 'Biologie', 'Biologie'),
    @(129, 'KGEO', 'SZZ - Geografie', 'Geografie'),
    @(130, 'KI', 'SZZ - Informatika', 'Informatika'),
    @(131, 'KCH', 'Chemie', 'Chemie'),
    @(132, 'KMA', 'Matematika', 'Matematika'),
    @(133, 'KMA', 'Matematická informatika', 'Matematická informatika'),
    @(134, 'KGEO', 'SZZ - Obecná geografie', 'Obecná geografie'),
    @(135, 'KFY', 'Elektronika a elektrotechnika', 'Elektronika a elektrotechnika'),
    @(140, 'KCH', 'Chemie', 'Chemie'),
    @(141, 'KFY', 'Fyzika', 'Fyzika'),
    @(142, 'KCH', 'Toxikologie', 'Toxikologie'),
    @(143, 'KFY', 'Počítačové modelování', 'Počítačové modelování'),
    @(150, 'KBI', 'Fyziologie, biochemie a mol. bio. buňky', 'Fyziologie, biochemie a molekulární biologie buňky'),
    @(151, 'KFY', 'Počítačové modelování', 'Počítačové modelování'),
    @(159, 'KBI', 'Učitelství biologie pro střední školy', 'Učitelství biologie pro střední školy'),
    @(160, 'KFY', 'Učitelství fyziky pro střední školy', 'Učitelství fyziky pro střední školy'),
    @(161, 'KCH', 'Učitelství chemie pro střední školy', 'Učitelství chemie pro střední školy'),
    @(163, 'KI', 'Učitelství informatiky pro střední školy', 'Učitelství informatiky pro střední školy'),
    @(164, 'KMA', 'Učitelství matematiky pro střední školy', 'Učitelství matematiky pro střední školy')
)

foreach ($chg in $changes) {
    $r = $chg[0]
    $ws.Cells.Item($r, 1).Value = $chg[1]
    $ws.Cells.Item($r, 4).Value = $chg[2]
    $ws.Cells.Item($r, 5).Value = $chg[3]
}
